$wb = $excel.ActiveWorkbook

# --- 1. Q9_Wagner sheet: select the used range (rows 1-3, all columns) ---
$q9Wagner = $wb.Worksheets.Item("Q9_Wagner")
$null = $q9Wagner.Range("A1:XFD3").Select()

# --- 2. Alex sheet: copy row 5 down into a new row 6 ---
$alex = $wb.Worksheets.Item("Alex")
$alex.Range("A6").Value = "dear god"
$alex.Range("B6").Value = "please work"
$null = $alex.Range("B5:B6").Select()

# --- 3. Insert a new sheet "Q9_Alex" right after "Alex" (mirrors Q9_Wagner) ---
$q9Alex = $wb.Worksheets.Add($null, $alex)
$q9Alex.Name = "Q9_Alex"

$q9Alex.Range("C1").Value = "Q9.1"
$q9Alex.Range("D1").Value = "Q9.2"
$q9Alex.Range("E1").Value = "Q9.3"
$q9Alex.Range("F1").Value = "Q9.4"
$q9Alex.Range("G1").Value = "Q9.5"

$q9Alex.Range("A2").Value = "Your name"
$q9Alex.Range("C2").Value = "hours"
$q9Alex.Range("E2").Value = "hours"

$q9Alex.Range("A3").Value = "Inputer"
$q9Alex.Range("B3").Value = "ID_fisher"
$q9Alex.Range("C3").Value = "soak_time_per_throw"
$q9Alex.Range("D3").Value = "num_throws"
$q9Alex.Range("E3").Value = "total_soak"
$q9Alex.Range("F3").Value = "gear"
$q9Alex.Range("G3").Value = "species"

# Match the bold "ID_fisher" header style used on the Q9_Wagner sheet
$q9Wagner.Range("B3").Copy()
$null = $q9Alex.Range("B3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Select the whole used range on the new sheet - it becomes the active sheet/tab
$null = $q9Alex.Range("A1:XFD3").Select()

# --- 4. Diana keeps its own selection/content; it's simply no longer the active tab,
#        which happened automatically once Q9_Alex was activated above. ---
